$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.614.06'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '2.291.29'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.50'
$ws.Range("E5").Value = '  +4.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '268.22'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.82'
$ws.Range("E10").Value = '  +2.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.91'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("D14").Value = '2.633.58'
$ws.Range("E14").Value = '  -1.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.19'
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.850'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = '2.289.30'
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").Value = '43.593.44'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.24'
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.53'
$ws.Range("E22").Value = '  +11.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.86'
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.11'
$ws.Range("E24").Value = '  -6.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.55'
$ws.Range("E26").Value = '  +2.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.21'
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("E28").Value = '  +2.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.03'
$ws.Range("E29").Value = '  +2.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  -4.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.47'
$ws.Range("E31").Value = '  +2.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.25'
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0894'
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.36'
$ws.Range("E34").Value = '  -3.05%  '
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0353'
$ws.Range("E37").Value = '  +1.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.37'
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.39'
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.244'
$ws.Range("E40").Value = '  +2.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.32'
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("E42").Value = '  +3.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.31'
$ws.Range("E43").Value = '  +1.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.03'
$ws.Range("E44").Value = '  +6.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.80'
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("E46").Value = '  -4.21%  '
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.30'
$ws.Range("E48").Value = '  -2.85%  '
$ws.Range("E49").Value = '  +0.21%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.50'
$ws.Range("E50").Value = '  +10.00%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.185'
$ws.Range("E51").Value = '  +9.06%  '

Write-Host "Applied all changes"